$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-19 Wednesday", "2025-02-20 Thursday"),
    @("652÷7=", "445÷2="),
    @("949÷4=", "779÷7="),
    @("472÷9=", "132÷9="),
    @("695÷7=", "798÷7="),
    @("999÷5=", "640÷3="),
    @("514÷2=", "297÷4="),
    @("630÷5=", "523÷7="),
    @("481÷7=", "363÷3="),
    @("464÷3=", "303÷7="),
    @("117÷7=", "591÷4="),
    @("769÷5=", "457÷6="),
    @("711÷2=", "116÷2="),
    @("591÷8=", "696÷5="),
    @("895÷5=", "432÷7="),
    @("884÷9=", "883÷4="),
    @("838÷5=", "430÷2="),
    @("940÷2=", "870÷7="),
    @("900÷7=", "944÷6="),
    @("696÷7=", "449÷9="),
    @("971÷8=", "439÷9="),
    @("284÷7=", "252÷5="),
    @("965÷3=", "165÷2="),
    @("679÷5=", "415÷2="),
    @("900÷2=", "305÷3="),
    @("901÷8=", "759÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
